# Weekly update: insert two new price records (new sampling date 45246)
# at the top of the data table (rows 17-18), pushing the existing rows
# 17-39 down to 19-41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 17 (shifts rows 17:39 -> 19:41,
# carrying their existing values/formatting down with them).
$ws.Rows("17:18").Insert()

# Common / constant column values shared by every data row in this table.
$mercadoId   = 7
$mercado     = "Terminal Hortofrutícola Agro Chillán"
$region      = "Ñuble"
$codreg      = 16
$categoriaId = 100112039
$categoria   = "Ciboulette"
$variedad    = "Sin especificar"
$unidad      = "`$/docena de atados"
$origen      = "Región Metropolitana"
$kgUnidades  = 3
$clasif      = "Hortaliza"

# New row 17: Primera, volumen 300, precios 2500/2500/2500, $/Kg 833
$ws.Cells.Item(17, 1).Value  = $mercadoId
$ws.Cells.Item(17, 2).Value  = $mercado
$ws.Cells.Item(17, 3).Value  = $region
$ws.Cells.Item(17, 4).Value  = 45246
$ws.Cells.Item(17, 5).Value  = $codreg
$ws.Cells.Item(17, 6).Value  = $categoriaId
$ws.Cells.Item(17, 7).Value  = $categoria
$ws.Cells.Item(17, 8).Value  = $variedad
$ws.Cells.Item(17, 9).Value  = "Primera"
$ws.Cells.Item(17, 10).Value = 300
$ws.Cells.Item(17, 11).Value = 2500
$ws.Cells.Item(17, 12).Value = 2500
$ws.Cells.Item(17, 13).Value = 2500
$ws.Cells.Item(17, 14).Value = $unidad
$ws.Cells.Item(17, 15).Value = $origen
$ws.Cells.Item(17, 16).Value = 833
$ws.Cells.Item(17, 17).Value = $kgUnidades
$ws.Cells.Item(17, 18).Value = $clasif

# New row 18: Segunda, volumen 200, precios 2000/2000/2000, $/Kg 667
$ws.Cells.Item(18, 1).Value  = $mercadoId
$ws.Cells.Item(18, 2).Value  = $mercado
$ws.Cells.Item(18, 3).Value  = $region
$ws.Cells.Item(18, 4).Value  = 45246
$ws.Cells.Item(18, 5).Value  = $codreg
$ws.Cells.Item(18, 6).Value  = $categoriaId
$ws.Cells.Item(18, 7).Value  = $categoria
$ws.Cells.Item(18, 8).Value  = $variedad
$ws.Cells.Item(18, 9).Value  = "Segunda"
$ws.Cells.Item(18, 10).Value = 200
$ws.Cells.Item(18, 11).Value = 2000
$ws.Cells.Item(18, 12).Value = 2000
$ws.Cells.Item(18, 13).Value = 2000
$ws.Cells.Item(18, 14).Value = $unidad
$ws.Cells.Item(18, 15).Value = $origen
$ws.Cells.Item(18, 16).Value = 667
$ws.Cells.Item(18, 17).Value = $kgUnidades
$ws.Cells.Item(18, 18).Value = $clasif
